# Applies the commit:
#   - rewrites the first sentence of the last content paragraph
#   - splits that paragraph in two, right after "...profesor." (keeping
#     the _GoBack bookmark attached to the end of the first half)
#   - the old second run's text becomes the start of the new paragraph,
#     prefixed with "Desafortunadamente"

$d = $word.ActiveDocument

# --- Step 0: the existing "_GoBack" bookmark sits right where we are about
# to split the paragraph; pull it out now and re-create it in the right
# spot once the surrounding text/paragraph edits are done. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 1: "Desafortunadamente" -> the new long sentence ---
$f1 = $d.Content
$f1.Find.Execute("Desafortunadamente", $true, $false, $false, $false, $false, $true, 1, $false, `
    "La información que el usuario puede exportar es el expediente académico del alumno y la lista de autorizados de una excursión de un profesor.", `
    2) | Out-Null

# --- Step 2: ", no funciona..." -> "Desafortunadamente, no funciona..." ---
$f2 = $d.Content
$f2.Find.Execute(", no funciona en el navegador Edge, sin embargo funciona correctamente en los navegadores Firefox y Chrome.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Desafortunadamente, no funciona en el navegador Edge, sin embargo funciona correctamente en los navegadores Firefox y Chrome.", `
    2) | Out-Null

# --- Step 3: split the paragraph right after "...profesor." ---
$fb = $d.Content
$fb.Find.Execute("profesor.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $fb.End
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# --- Step 4: re-insert "_GoBack" at the end of the first paragraph (right
# after "...profesor.", before the new paragraph mark). Directly adding a
# zero-length bookmark exactly at a paragraph boundary is unreliable, so a
# throwaway marker run is inserted, the bookmark is wrapped tightly around
# it (a safe, non-boundary range), and then the marker text is deleted --
# the bookmark collapses to the correct zero-width position in the process. ---
$markerRange = $d.Range($splitPos, $splitPos)
$markerRange.InsertBefore("@@MARK@@")

$fm = $d.Content
$fm.Find.Execute("@@MARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $fm) | Out-Null

$fm2 = $d.Content
$fm2.Find.Execute("@@MARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
